# Weekly driver report update for 2025-04-20
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bad Drivers table: refreshed "Good Roaming Calculation (%)" figures ---
$ws.Range("D3").Value = 88.40000000000001
$ws.Range("D4").Value = 97.90000000000001

# --- Good Drivers table: new week's data, one new driver inserted at the
#     top (row 13) and every other row shifted down by one, with each
#     row's "Driver Vintage" date now populated. ---

# Row 13 (new entry)
$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B13").Value = 445055
$ws.Range("D13").Value = 99.90000000000001
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2024-11-10"

# Row 14
$ws.Range("A14").Value = "MediaTek MT7921 Wi-Fi 6 802.11ax PCIe Adapter - 3.0.1.1255"
$ws.Range("B14").Value = 23159
$ws.Range("D14").Value = 99.90000000000001
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2022-07-04"

# Row 15
$ws.Range("A15").Value = "MediaTek MT7921 Wi-Fi 6 802.11ax PCIe Adapter - 3.0.1.1216"
$ws.Range("B15").Value = 36106
$ws.Range("D15").Value = 100
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2021-12-23"

# Row 16
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B16").Value = 77849
$ws.Range("D16").Value = 99.90000000000001
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2021-08-18"

# Row 17
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B17").Value = 34244
$ws.Range("D17").Value = 100
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2021-04-27"

# Row 18 is unchanged (Intel 21.110.3.2 stays in place)

# Row 19
$ws.Range("A19").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B19").Value = 113652
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2020-01-06"

# Row 20
$ws.Range("A20").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B20").Value = 56018
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2019-12-14"
